$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect the latest scrape.
# Rows 21 and 22 swapped ranking order (WrappedliquidstakedEther2.0 now
# ranks above Dai), and associated price/link/volume values updated.
# Numeric-looking Price values are written as text (matching the workbook
# convention of storing Price as inline-string) by forcing a Text number
# format before the assignment, then restoring the default cell style so
# no visible formatting changes are introduced.

$ws.Range("D2").Value = "29.368.61"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.872.77"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7163"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07798"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3068"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08250"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.877.67"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7231"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "29.458.43"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.845"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007864"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.138.40"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.771"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1558"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.929"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.357"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.479"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.334"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.090"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05247"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7181"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.680"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01866"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.710"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").Value = "1.183.68"
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9106"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.021"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4301"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5366"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.764"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.161"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.024"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.64%  "
